$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 42266
$ws.Range("F4").Value = 9976
$ws.Range("F5").Value = 223
$ws.Range("F6").Value = 1032
$ws.Range("F8").Value = 774
$ws.Range("F13").Value = 134
$ws.Range("F15").Value = 347
$ws.Range("F16").Value = 1604
$ws.Range("F18").Value = 785
$ws.Range("F19").Value = 756
$ws.Range("F20").Value = 498
$ws.Range("F26").Value = 570
$ws.Range("F27").Value = 565
$ws.Range("F32").Value = 455
$ws.Range("F34").Value = 231
$ws.Range("F37").Value = 1415
$ws.Range("F39").Value = 1311

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 470

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 9976
$ws.Range("F8").Value = 223
$ws.Range("F9").Value = 1032
$ws.Range("F10").Value = 1032
$ws.Range("F18").Value = 134
$ws.Range("F20").Value = 347
$ws.Range("F21").Value = 1604
$ws.Range("F23").Value = 756
$ws.Range("F24").Value = 498
$ws.Range("F29").Value = 570
$ws.Range("F32").Value = 565
$ws.Range("F37").Value = 455
$ws.Range("F39").Value = 231
$ws.Range("F41").Value = 1311
